# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage, so values that look like
# plain numbers (e.g. "44.61") are not auto-converted to numeric cells by Excel.
# The NumberFormat is reset back to the default afterwards so no stray number
# format is left applied to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# cell reference -> new value (column D numeric-looking values are flagged so
# they get routed through Set-TextValue).
$updates = @(
    @('D2', '28.410.15', $false),
    @('E2', '  -0.25%  ', $false),
    @('D3', '1.572.86', $false),
    @('E3', '  +0.06%  ', $false),
    @('E4', '  -0.22%  ', $false),
    @('D5', '211.88', $true),
    @('E5', '  -0.06%  ', $false),
    @('D6', '0.492', $true),
    @('E6', '  -0.10%  ', $false),
    @('E7', '  -0.20%  ', $false),
    @('D8', '44.61', $true),
    @('D9', '23.70', $true),
    @('E9', '  -1.73%  ', $false),
    @('E10', '  -0.56%  ', $false),
    @('E11', '  -0.59%  ', $false),
    @('E12', '  +1.38%  ', $false),
    @('D13', '1.798.35', $false),
    @('E13', '  +0.01%  ', $false),
    @('D14', '1.570.05', $false),
    @('E14', '  -0.18%  ', $false),
    @('E15', '  -0.26%  ', $false),
    @('D16', '28.414.96', $false),
    @('E16', '  -0.32%  ', $false),
    @('D17', '0.515', $true),
    @('E17', '  -1.04%  ', $false),
    @('D18', '61.63', $true),
    @('E18', '  -0.80%  ', $false),
    @('D19', '229.97', $true),
    @('E19', '  +1.13%  ', $false),
    @('E20', '  +0.40%  ', $false),
    @('D21', '0.0₃0683', $false),
    @('E21', '  -1.28%  ', $false),
    @('E22', '  -0.16%  ', $false),
    @('D23', '3.96', $true),
    @('E23', '  +1.93%  ', $false),
    @('E24', '  -1.08%  ', $false),
    @('E25', '  +1.61%  ', $false),
    @('D26', '151.35', $true),
    @('E26', '  +0.08%  ', $false),
    @('D27', '14.92', $true),
    @('E27', '  -0.29%  ', $false),
    @('E28', '  -0.32%  ', $false),
    @('E29', '  -1.30%  ', $false),
    @('E30', '  -0.20%  ', $false),
    @('D31', '0.0483', $true),
    @('E31', '  +4.26%  ', $false),
    @('E32', '  -3.09%  ', $false),
    @('E33', '  -0.59%  ', $false),
    @('D34', '3.10', $true),
    @('E34', '  -1.24%  ', $false),
    @('D35', '1.381.12', $false),
    @('E35', '  -0.70%  ', $false),
    @('E36', '  +4.38%  ', $false),
    @('D37', '1.51', $true),
    @('E37', '  -1.91%  ', $false),
    @('E38', '  +0.09%  ', $false),
    @('E39', '  +1.34%  ', $false),
    @('E40', '  -1.59%  ', $false),
    @('D41', '0.522', $true),
    @('E41', '  -1.93%  ', $false),
    @('B42', 'PaxDollar', $false),
    @('C42', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', $false),
    @('D42', '1.00', $true),
    @('E42', '  -0.20%  ', $false),
    @('B43', 'RenderToken', $false),
    @('C43', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', $false),
    @('D43', '1.90', $true),
    @('E43', '  +2.60%  ', $false),
    @('D44', '0.787', $true),
    @('E44', '  -0.75%  ', $false),
    @('E45', '  +1.39%  ', $false),
    @('E46', '  -4.44%  ', $false),
    @('D47', '62.35', $true),
    @('E47', '  -1.02%  ', $false),
    @('E48', '  -6.20%  ', $false),
    @('D49', '1.710.43', $false),
    @('E49', '  -0.02%  ', $false),
    @('E50', '  -0.37%  ', $false),
    @('D51', '85.32', $true),
    @('E51', '  -0.83%  ', $false)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $value = $u[1]
    $forceText = $u[2]
    $range = $ws.Range($ref)
    if ($forceText) {
        Set-TextValue $range $value
    } else {
        $range.Value = $value
    }
}